$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting existing rows 48:164 down to 49:165
$ws.Rows("48:48").Insert()

# Fill in the new row 48 with the new record's data
$ws.Range("A48").Value = 4
$ws.Range("B48").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C48").Value = "Los Lagos"
$ws.Range("D48").Value = 44519
$ws.Range("D48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E48").Value = 10
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100108
$ws.Range("H48").Value = "Tropicales y subtropicales"
$ws.Range("I48").Value = 100108005
$ws.Range("J48").Value = "Piña"
$ws.Range("K48").Value = "Caramelo"
$ws.Range("L48").Value = "Tercera"
$ws.Range("M48").Value = 200
$ws.Range("N48").Value = 20000
$ws.Range("O48").Value = 21000
$ws.Range("P48").Value = 20500
$ws.Range("Q48").Value = "$/caja 16 unidades"
$ws.Range("R48").Value = "Ecuador"
$ws.Range("S48").Value = 1281
$ws.Range("T48").Value = 16
